# Apply cryptocurrency price/volume updates as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'68.856.59"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.67%  '
$ws.Range('D3').Value = "'3.716.21"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.25%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'605.03"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').Value = "'182.60"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.55%  '
$ws.Range('D7').Value = "'3.714.05"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.22%  '
$ws.Range('E8').Value = '  -5.61%  '
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('E10').Value = '  -3.89%  '
$ws.Range('E11').Value = '  -8.71%  '
$ws.Range('D12').Value = "'56.94"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.54%  '
$ws.Range('D13').Value = "'0.0000295"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -8.75%  '
$ws.Range('D14').Value = "'10.69"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.08%  '
$ws.Range('D15').Value = "'4.317.09"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.00%  '
$ws.Range('D16').Value = "'3.719.32"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.04%  '
$ws.Range('D17').Value = "'19.50"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.58%  '
$ws.Range('E18').Value = '  -1.98%  '
$ws.Range('D19').Value = "'12.98"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.03%  '
$ws.Range('E20').Value = '  -7.13%  '
$ws.Range('D21').Value = "'68.775.97"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.52%  '
$ws.Range('D22').Value = "'415.16"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.85%  '
$ws.Range('D23').Value = "'4.67"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.41%  '
$ws.Range('D24').Value = "'89.17"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.69%  '
$ws.Range('D25').Value = "'3.06"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.79%  '
$ws.Range('D26').Value = "'12.75"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -8.31%  '
$ws.Range('D27').Value = "'10.97"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.31%  '
$ws.Range('E28').Value = '  -3.66%  '
$ws.Range('D29').Value = "'6.07"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.79%  '
$ws.Range('D30').Value = "'9.64"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.44%  '
$ws.Range('E31').Value = '  -6.09%  '
$ws.Range('D32').Value = "'7.33"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -15.20%  '
$ws.Range('E33').Value = '  -7.77%  '
$ws.Range('E34').Value = '  -5.31%  '
$ws.Range('D35').Value = "'44.00"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.88%  '
$ws.Range('D36').Value = "'65.09"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.76%  '
$ws.Range('D37').Value = "'605.42"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.34%  '
$ws.Range('D38').Value = "'0.0₃0888"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -11.29%  '
$ws.Range('D39').Value = "'0.407"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.42%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = "'1.00"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('E42').Value = '  -5.69%  '
$ws.Range('D43').Value = "'3.07"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.68%  '
$ws.Range('D44').Value = "'0.0443"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.31%  '
$ws.Range('D45').Value = "'2.67"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.12%  '
$ws.Range('D46').Value = "'2.81"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -10.85%  '
$ws.Range('D47').Value = "'9.25"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.46%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = "'0.136"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.25%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').Value = "'2.74"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.89%  '
$ws.Range('D50').Value = "'2.780.60"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.27%  '
$ws.Range('D51').Value = "'3.09"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.39%  '
